# ------------------------------------------------------------------
# Applies the "Add files via upload" edit: the last duplicate test-step
# row (testTp4, step 2) on the Login, Logout and Home sheets is turned
# into a "testTp5" row (step 1) and five brand-new rows are appended
# for testTp6..testTp10 (also step 1). Finally the Logout sheet becomes
# the active sheet/tab.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------- Login sheet ----------------
$login = $wb.Worksheets.Item("Login")

# New row for testTp6 (copy formatting from the row we are about to edit)
$login.Rows(16).Copy() | Out-Null
$login.Rows(17).Insert() | Out-Null
$login.Cells.Item(17,1).Value = "testTp6"
$login.Cells.Item(17,2).Value = 1
$login.Cells.Item(17,3).Value = 1
$login.Cells.Item(17,4).Value = "HLWorkerContraCosta"

# Turn the old duplicate "testTp4" (step 2) row into "testTp5" (step 1)
$login.Cells.Item(16,1).Value = "testTp5"
$login.Cells.Item(16,2).Value = 1
$login.Cells.Item(16,4).Value = "HLWorkerContraCosta"

# New rows for testTp7 .. testTp10
$login.Rows(18).Insert() | Out-Null
$login.Cells.Item(18,1).Value = "testTp7"
$login.Cells.Item(18,2).Value = 1
$login.Cells.Item(18,3).Value = 1
$login.Cells.Item(18,4).Value = "HLWorkerContraCosta"

$login.Rows(19).Insert() | Out-Null
$login.Cells.Item(19,1).Value = "testTp8"
$login.Cells.Item(19,2).Value = 1
$login.Cells.Item(19,3).Value = 1
$login.Cells.Item(19,4).Value = "HLWorkerContraCosta"

$login.Rows(20).Insert() | Out-Null
$login.Cells.Item(20,1).Value = "testTp9"
$login.Cells.Item(20,2).Value = 1
$login.Cells.Item(20,3).Value = 1
$login.Cells.Item(20,4).Value = "HLWorkerContraCosta"

$login.Rows(21).Insert() | Out-Null
$login.Cells.Item(21,1).Value = "testTp10"
$login.Cells.Item(21,2).Value = 1
$login.Cells.Item(21,3).Value = 1
$login.Cells.Item(21,4).Value = "HLWorkerContraCosta"

$login.Range("A16:A21").Select() | Out-Null

# ---------------- Home sheet ----------------
$home = $wb.Worksheets.Item("Home")

$home.Rows(11).Copy() | Out-Null
$home.Rows(12).Insert() | Out-Null
$home.Cells.Item(12,1).Value = "testTp6"
$home.Cells.Item(12,2).Value = 1
$home.Cells.Item(12,3).Value = 1
$home.Cells.Item(12,4).Value = "Click"
$home.Cells.Item(12,5).Value = "Click"

$home.Cells.Item(11,1).Value = "testTp5"
$home.Cells.Item(11,2).Value = 1

$home.Rows(13).Insert() | Out-Null
$home.Cells.Item(13,1).Value = "testTp7"
$home.Cells.Item(13,2).Value = 1
$home.Cells.Item(13,3).Value = 1
$home.Cells.Item(13,4).Value = "Click"
$home.Cells.Item(13,5).Value = "Click"

$home.Rows(14).Insert() | Out-Null
$home.Cells.Item(14,1).Value = "testTp8"
$home.Cells.Item(14,2).Value = 1
$home.Cells.Item(14,3).Value = 1
$home.Cells.Item(14,4).Value = "Click"
$home.Cells.Item(14,5).Value = "Click"

$home.Rows(15).Insert() | Out-Null
$home.Cells.Item(15,1).Value = "testTp9"
$home.Cells.Item(15,2).Value = 1
$home.Cells.Item(15,3).Value = 1
$home.Cells.Item(15,4).Value = "Click"
$home.Cells.Item(15,5).Value = "Click"

$home.Rows(16).Insert() | Out-Null
$home.Cells.Item(16,1).Value = "testTp10"
$home.Cells.Item(16,2).Value = 1
$home.Cells.Item(16,3).Value = 1
$home.Cells.Item(16,4).Value = "Click"
$home.Cells.Item(16,5).Value = "Click"

$home.Range("A11:A16").Select() | Out-Null

# ---------------- Logout sheet ----------------
$logout = $wb.Worksheets.Item("Logout")

# Row 9 (not 10) is the formatting source: column F on row 10 lost its
# explicit style once the sheet was re-saved, matching row 9's shape.
$logout.Rows(9).Copy() | Out-Null
$logout.Rows(10).Insert() | Out-Null
$logout.Cells.Item(10,1).Value = "testTp6"
$logout.Cells.Item(10,2).Value = 1
$logout.Cells.Item(10,3).Value = 1
$logout.Cells.Item(10,4).Value = "Click"
$logout.Cells.Item(10,5).Value = "Click"
$logout.Cells.Item(10,6).Value = "Yes"

$logout.Cells.Item(11,1).Value = "testTp5"
$logout.Cells.Item(11,2).Value = 1

$logout.Rows(12).Insert() | Out-Null
$logout.Cells.Item(12,1).Value = "testTp7"
$logout.Cells.Item(12,2).Value = 1
$logout.Cells.Item(12,3).Value = 1
$logout.Cells.Item(12,4).Value = "Click"
$logout.Cells.Item(12,5).Value = "Click"
$logout.Cells.Item(12,6).Value = "Yes"

$logout.Rows(13).Insert() | Out-Null
$logout.Cells.Item(13,1).Value = "testTp8"
$logout.Cells.Item(13,2).Value = 1
$logout.Cells.Item(13,3).Value = 1
$logout.Cells.Item(13,4).Value = "Click"
$logout.Cells.Item(13,5).Value = "Click"
$logout.Cells.Item(13,6).Value = "Yes"

$logout.Rows(14).Insert() | Out-Null
$logout.Cells.Item(14,1).Value = "testTp9"
$logout.Cells.Item(14,2).Value = 1
$logout.Cells.Item(14,3).Value = 1
$logout.Cells.Item(14,4).Value = "Click"
$logout.Cells.Item(14,5).Value = "Click"
$logout.Cells.Item(14,6).Value = "Yes"

$logout.Rows(15).Insert() | Out-Null
$logout.Cells.Item(15,1).Value = "testTp10"
$logout.Cells.Item(15,2).Value = 1
$logout.Cells.Item(15,3).Value = 1
$logout.Cells.Item(15,4).Value = "Click"
$logout.Cells.Item(15,5).Value = "Click"
$logout.Cells.Item(15,6).Value = "Yes"

$logout.Range("F10").Select() | Out-Null
